$wb = $excel.ActiveWorkbook

# A pristine cell (sheet1 K2) already holds the literal text "2013-12-03".
# Copying it (instead of assigning the literal string through .Value) avoids
# the automatic "looks like a date" conversion that the COM layer otherwise
# performs on ISO-formatted date strings.
$dateDonor = $wb.Worksheets.Item(1).Cells.Item(2,11)

# ===================================================================
# Sheet 6 - insurance (保險)
#   existing: A=index B=company C=name D=owner
#   new:      E=property_category F=category G=date H=legislator_name
#             I=legislator_id J=source_file K=index
# ===================================================================
$ws6 = $wb.Worksheets.Item(6)

# Copy header style (bold + border, same as B1) onto the new header cells
# BEFORE touching B1's own value/content.
$ws6.Cells.Item(1,2).Copy($ws6.Cells.Item(1,5))
$ws6.Cells.Item(1,2).Copy($ws6.Cells.Item(1,6))
$ws6.Cells.Item(1,2).Copy($ws6.Cells.Item(1,7))
$ws6.Cells.Item(1,2).Copy($ws6.Cells.Item(1,8))
$ws6.Cells.Item(1,2).Copy($ws6.Cells.Item(1,9))
$ws6.Cells.Item(1,2).Copy($ws6.Cells.Item(1,10))
$ws6.Cells.Item(1,2).Copy($ws6.Cells.Item(1,11))

# header row 1 labels
$ws6.Cells.Item(1,2).Value = "company"
$ws6.Cells.Item(1,3).Value = "name"
$ws6.Cells.Item(1,4).Value = "owner"
$ws6.Cells.Item(1,5).Value = "property_category"
$ws6.Cells.Item(1,6).Value = "category"
$ws6.Cells.Item(1,7).Value = "date"
$ws6.Cells.Item(1,8).Value = "legislator_name"
$ws6.Cells.Item(1,9).Value = "legislator_id"
$ws6.Cells.Item(1,10).Value = "source_file"
$ws6.Cells.Item(1,11).Value = "index"

# data rows 2-8
$ins_company = @("台灣人壽","台灣人壽","國泰人壽","國泰人壽","國泰人壽","國泰人壽","國泰人壽")
$ins_name    = @("歲歲長泰還本終身險","新祥和定期壽險","得意還本終身險","富貴年年終身險","添寶養老壽險","雙好還本終身險","雙星還本終身險")
$ins_owner   = @("李永得","李永得","邱議瑩","邱議瑩","邱議瑩","邱議瑩","邱議瑩")
$ins_index   = @(97,98,99,100,101,102,103)

for ($k = 0; $k -lt 7; $k++) {
    $r = $k + 2
    $ws6.Cells.Item($r,2).Value = $ins_company[$k]
    $ws6.Cells.Item($r,3).Value = $ins_name[$k]
    $ws6.Cells.Item($r,4).Value = $ins_owner[$k]
    $ws6.Cells.Item($r,5).Value = "insurance"
    $ws6.Cells.Item($r,6).Value = "normal"
    $dateDonor.Copy($ws6.Cells.Item($r,7))
    $ws6.Cells.Item($r,8).Value = "邱議瑩"
    $ws6.Cells.Item($r,9).Value = 913
    $ws6.Cells.Item($r,10).Value = "tmp40191"
    $ws6.Cells.Item($r,11).Value = $ins_index[$k]
}

# ===================================================================
# Sheet 7 - debt (債務)
#   existing: A=index B=species C=debtor D=owner E=total
#             F=register_date G=register_reason
#   new:      H=property_category I=category J=date K=legislator_name
#             L=legislator_id M=source_file N=index
# ===================================================================
$ws7 = $wb.Worksheets.Item(7)

$ws7.Cells.Item(1,2).Copy($ws7.Cells.Item(1,8))
$ws7.Cells.Item(1,2).Copy($ws7.Cells.Item(1,9))
$ws7.Cells.Item(1,2).Copy($ws7.Cells.Item(1,10))
$ws7.Cells.Item(1,2).Copy($ws7.Cells.Item(1,11))
$ws7.Cells.Item(1,2).Copy($ws7.Cells.Item(1,12))
$ws7.Cells.Item(1,2).Copy($ws7.Cells.Item(1,13))
$ws7.Cells.Item(1,2).Copy($ws7.Cells.Item(1,14))

$ws7.Cells.Item(1,2).Value = "species"
$ws7.Cells.Item(1,3).Value = "debtor"
$ws7.Cells.Item(1,4).Value = "owner"
$ws7.Cells.Item(1,5).Value = "total"
$ws7.Cells.Item(1,6).Value = "register_date"
$ws7.Cells.Item(1,7).Value = "register_reason"
$ws7.Cells.Item(1,8).Value = "property_category"
$ws7.Cells.Item(1,9).Value = "category"
$ws7.Cells.Item(1,10).Value = "date"
$ws7.Cells.Item(1,11).Value = "legislator_name"
$ws7.Cells.Item(1,12).Value = "legislator_id"
$ws7.Cells.Item(1,13).Value = "source_file"
$ws7.Cells.Item(1,14).Value = "index"

# data rows 2-3 (B:G unchanged values, just add H:N)
$debt_index = @(116,117)
for ($k = 0; $k -lt 2; $k++) {
    $r = $k + 2
    $ws7.Cells.Item($r,8).Value = "debt"
    $ws7.Cells.Item($r,9).Value = "normal"
    $dateDonor.Copy($ws7.Cells.Item($r,10))
    $ws7.Cells.Item($r,11).Value = "邱議瑩"
    $ws7.Cells.Item($r,12).Value = 913
    $ws7.Cells.Item($r,13).Value = "tmp40191"
    $ws7.Cells.Item($r,14).Value = $debt_index[$k]
}

# ===================================================================
# Sheet 8 - investment (事業投資)
#   existing: A=index B=owner C=company D=address E=total
#             F=register_date G=register_reason
#   new:      H=property_category I=category J=date K=legislator_name
#             L=legislator_id M=source_file N=index
# ===================================================================
$ws8 = $wb.Worksheets.Item(8)

$ws8.Cells.Item(1,2).Copy($ws8.Cells.Item(1,8))
$ws8.Cells.Item(1,2).Copy($ws8.Cells.Item(1,9))
$ws8.Cells.Item(1,2).Copy($ws8.Cells.Item(1,10))
$ws8.Cells.Item(1,2).Copy($ws8.Cells.Item(1,11))
$ws8.Cells.Item(1,2).Copy($ws8.Cells.Item(1,12))
$ws8.Cells.Item(1,2).Copy($ws8.Cells.Item(1,13))
$ws8.Cells.Item(1,2).Copy($ws8.Cells.Item(1,14))

$ws8.Cells.Item(1,2).Value = "owner"
$ws8.Cells.Item(1,3).Value = "company"
$ws8.Cells.Item(1,4).Value = "address"
$ws8.Cells.Item(1,5).Value = "total"
$ws8.Cells.Item(1,6).Value = "register_date"
$ws8.Cells.Item(1,7).Value = "register_reason"
$ws8.Cells.Item(1,8).Value = "property_category"
$ws8.Cells.Item(1,9).Value = "category"
$ws8.Cells.Item(1,10).Value = "date"
$ws8.Cells.Item(1,11).Value = "legislator_name"
$ws8.Cells.Item(1,12).Value = "legislator_id"
$ws8.Cells.Item(1,13).Value = "source_file"
$ws8.Cells.Item(1,14).Value = "index"

# data rows 2-3 (B:G unchanged values, just add H:N)
$inv_index = @(122,123)
for ($k = 0; $k -lt 2; $k++) {
    $r = $k + 2
    $ws8.Cells.Item($r,8).Value = "investment"
    $ws8.Cells.Item($r,9).Value = "normal"
    $dateDonor.Copy($ws8.Cells.Item($r,10))
    $ws8.Cells.Item($r,11).Value = "邱議瑩"
    $ws8.Cells.Item($r,12).Value = 913
    $ws8.Cells.Item($r,13).Value = "tmp40191"
    $ws8.Cells.Item($r,14).Value = $inv_index[$k]
}
